$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Updated data values in column C (Real GDP, billions Kc) for rows 2-22,
#    now formatted with 4 decimal places (numFmt "0.0000").
# ---------------------------------------------------------------------------
$cValues = @{
    2  = 2430.7012405999999
    3  = 2496.1811048
    4  = 2640.2801350999998
    5  = 2759.6210000000001
    6  = 2740.2440000000001
    7  = 2731.971
    8  = 2766.7359999999999
    9  = 2895.2339999999999
    10 = 2983.2489999999998
    11 = 3029.8110000000001
    12 = 3138.8580000000002
    13 = 3289.2240000000002
    14 = 3503.4989999999998
    15 = 3751.2109999999998
    16 = 3958.0729999999999
    17 = 4058.5740000000001
    18 = 3867.8029999999999
    19 = 3950.607
    20 = 4028.489
    21 = 3999.0720000000001
    22 = 3970.7130000000002
}

foreach ($row in $cValues.Keys) {
    $cell = $ws.Range("C$row")
    $cell.Value = $cValues[$row]
    $cell.NumberFormat = "0.0000"
}

# Row 23: column C previously empty, now gets a value (kept General format).
$ws.Range("C23").Value = 4056.26

# ---------------------------------------------------------------------------
# 2. Interest Rate column (Q) updates / new values.
# ---------------------------------------------------------------------------
$qValues = @{
    4  = 11
    5  = 11.3
    6  = 12.5
    7  = 14.9
    8  = 9.6
    9  = 5.43
    11 = 4.75
    12 = 2.5499999999999998
    13 = 2
    14 = 2.5
    15 = 2
    16 = 2.5
    17 = 3.5
    18 = 2.25
    19 = 1
    20 = 0.76
    21 = 0.76
    22 = 0.04
    23 = 0.04
}

foreach ($row in $qValues.Keys) {
    $ws.Range("Q$row").Value = $qValues[$row]
}

# Row 23: column R previously empty, now gets a value.
$ws.Range("R23").Value = -2.1

# ---------------------------------------------------------------------------
# 3. Row 8 gains a custom row height.
# ---------------------------------------------------------------------------
$ws.Rows("8").RowHeight = 15.75

# ---------------------------------------------------------------------------
# 4. New trailing rows (quiz-hint area) below the data table.
# ---------------------------------------------------------------------------
$ws.Rows("24").RowHeight = 20.25

$ws.Range("B26").NumberFormat = "#,##0"

$ws.Range("D27:D48").NumberFormat = "0.0000"

$boldRng = $ws.Range("C30:C48,D48")
$boldRng.Font.Name = "Arial CE"
$boldRng.Font.Size = 10
$boldRng.Font.Bold = $true
$ws.Range("C30:C48").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# 5. Selection moves to Q24 (matches the author's final cursor position).
# ---------------------------------------------------------------------------
$ws.Range("Q24").Select()
